$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles + row height) from the last existing data row (453)
# down through the new rows (454:489) before writing values.
$ws.Range("A453:F453").Copy()
$ws.Range("A454:F489").PasteSpecial(-4122)
$ws.Range("A454:F489").RowHeight = 13.55

# Widen column B to fit the longer program-name strings.
$ws.Columns.Item(2).ColumnWidth = 20

# Append the new export rows.
$ws.Cells.Item(454,1).Value = "2024-11-07 09:50:44"
$ws.Cells.Item(454,2).Value = "020071.NC"
$ws.Cells.Item(454,3).Value = 20
$ws.Cells.Item(454,4).Value = 12
$ws.Cells.Item(454,5).Value = 25
$ws.Cells.Item(454,6).Value = 1
$ws.Cells.Item(455,1).Value = "2024-11-11 15:54:01"
$ws.Cells.Item(455,2).Value = "020146.NC"
$ws.Cells.Item(455,3).Value = 1
$ws.Cells.Item(455,4).Value = 2
$ws.Cells.Item(455,5).Value = 12
$ws.Cells.Item(455,6).Value = 1
$ws.Cells.Item(456,1).Value = "2024-11-12 14:14:47"
$ws.Cells.Item(456,2).Value = "020090.slp"
$ws.Cells.Item(456,3).Value = 14
$ws.Cells.Item(456,4).Value = 208
$ws.Cells.Item(456,5).Value = 5
$ws.Cells.Item(456,6).Value = 1
$ws.Cells.Item(457,1).Value = "2024-11-13 11:37:59"
$ws.Cells.Item(457,2).Value = "020249.NC"
$ws.Cells.Item(457,3).Value = 24
$ws.Cells.Item(457,4).Value = 73
$ws.Cells.Item(457,5).Value = 10
$ws.Cells.Item(457,6).Value = 1
$ws.Cells.Item(458,1).Value = "2024-11-13 12:25:13"
$ws.Cells.Item(458,2).Value = "020275.NC"
$ws.Cells.Item(458,3).Value = 41
$ws.Cells.Item(458,4).Value = 34
$ws.Cells.Item(458,5).Value = 20
$ws.Cells.Item(458,6).Value = 1
$ws.Cells.Item(459,1).Value = "2024-11-13 13:30:47"
$ws.Cells.Item(459,2).Value = "020272.NC"
$ws.Cells.Item(459,3).Value = 13
$ws.Cells.Item(459,4).Value = 11
$ws.Cells.Item(459,5).Value = 20
$ws.Cells.Item(459,6).Value = 1
$ws.Cells.Item(460,1).Value = "2024-11-14 10:02:01"
$ws.Cells.Item(460,2).Value = "020314.NC"
$ws.Cells.Item(460,3).Value = 44
$ws.Cells.Item(460,4).Value = 264
$ws.Cells.Item(460,5).Value = 8
$ws.Cells.Item(460,6).Value = 1
$ws.Cells.Item(461,1).Value = "2024-11-14 10:14:23"
$ws.Cells.Item(461,2).Value = "020289.NC"
$ws.Cells.Item(461,3).Value = 30
$ws.Cells.Item(461,4).Value = 14
$ws.Cells.Item(461,5).Value = 2
$ws.Cells.Item(461,6).Value = 1
$ws.Cells.Item(462,1).Value = "2024-11-14 15:42:03"
$ws.Cells.Item(462,2).Value = "020360.NC"
$ws.Cells.Item(462,3).Value = 1
$ws.Cells.Item(462,4).Value = 5
$ws.Cells.Item(462,5).Value = 5
$ws.Cells.Item(462,6).Value = 1
$ws.Cells.Item(463,1).Value = "2024-11-14 15:52:34"
$ws.Cells.Item(463,2).Value = "020345.NC"
$ws.Cells.Item(463,3).Value = 47
$ws.Cells.Item(463,4).Value = 145
$ws.Cells.Item(463,5).Value = 8
$ws.Cells.Item(463,6).Value = 1
$ws.Cells.Item(464,1).Value = "2024-11-18 15:32:28"
$ws.Cells.Item(464,2).Value = "020521.NC"
$ws.Cells.Item(464,3).Value = 33
$ws.Cells.Item(464,4).Value = 18
$ws.Cells.Item(464,5).Value = 16
$ws.Cells.Item(464,6).Value = 1
$ws.Cells.Item(465,1).Value = "2024-11-19 13:26:14"
$ws.Cells.Item(465,2).Value = "020507.NC"
$ws.Cells.Item(465,3).Value = 71
$ws.Cells.Item(465,4).Value = 500
$ws.Cells.Item(465,5).Value = 6
$ws.Cells.Item(465,6).Value = 1
$ws.Cells.Item(466,1).Value = "2024-11-20 07:54:21"
$ws.Cells.Item(466,2).Value = "020496.NC"
$ws.Cells.Item(466,3).Value = 72
$ws.Cells.Item(466,4).Value = 500
$ws.Cells.Item(466,5).Value = 6
$ws.Cells.Item(466,6).Value = 1
$ws.Cells.Item(467,1).Value = "2024-11-20 11:43:15"
$ws.Cells.Item(467,2).Value = "020505.NC"
$ws.Cells.Item(467,3).Value = 72
$ws.Cells.Item(467,4).Value = 500
$ws.Cells.Item(467,5).Value = 6
$ws.Cells.Item(467,6).Value = 1
$ws.Cells.Item(468,1).Value = "2024-11-20 14:37:55"
$ws.Cells.Item(468,2).Value = "020503.NC"
$ws.Cells.Item(468,3).Value = 71
$ws.Cells.Item(468,4).Value = 500
$ws.Cells.Item(468,5).Value = 6
$ws.Cells.Item(468,6).Value = 1
$ws.Cells.Item(469,1).Value = "2024-11-07 11:06:01"
$ws.Cells.Item(469,2).Value = "020105.NC"
$ws.Cells.Item(469,3).Value = 1116
$ws.Cells.Item(469,4).Value = 824.4
$ws.Cells.Item(469,5).Value = 25
$ws.Cells.Item(469,6).Value = 1
$ws.Cells.Item(470,1).Value = "2024-11-20 17:27:26"
$ws.Cells.Item(470,2).Value = "020501.NC"
$ws.Cells.Item(470,3).Value = 72
$ws.Cells.Item(470,4).Value = 500
$ws.Cells.Item(470,5).Value = 6
$ws.Cells.Item(470,6).Value = 1
$ws.Cells.Item(471,1).Value = "2024-11-20 21:01:32"
$ws.Cells.Item(471,2).Value = "020908.NC"
$ws.Cells.Item(471,3).Value = 576
$ws.Cells.Item(471,4).Value = 158
$ws.Cells.Item(471,5).Value = 32
$ws.Cells.Item(471,6).Value = 1
$ws.Cells.Item(472,1).Value = "2024-11-21 19:57:12"
$ws.Cells.Item(472,2).Value = "020509.NC"
$ws.Cells.Item(472,3).Value = 86
$ws.Cells.Item(472,4).Value = 594
$ws.Cells.Item(472,5).Value = 6
$ws.Cells.Item(472,6).Value = 1
$ws.Cells.Item(473,1).Value = "2024-11-22 02:47:57"
$ws.Cells.Item(473,2).Value = "020490.NC"
$ws.Cells.Item(473,3).Value = 71
$ws.Cells.Item(473,4).Value = 500
$ws.Cells.Item(473,5).Value = 6
$ws.Cells.Item(473,6).Value = 1
$ws.Cells.Item(474,1).Value = "2024-11-25 17:07:02"
$ws.Cells.Item(474,2).Value = "020508.NC"
$ws.Cells.Item(474,3).Value = 78
$ws.Cells.Item(474,4).Value = 529
$ws.Cells.Item(474,5).Value = 6
$ws.Cells.Item(474,6).Value = 1
$ws.Cells.Item(475,1).Value = "2024-11-25 21:08:54"
$ws.Cells.Item(475,2).Value = "020494.NC"
$ws.Cells.Item(475,3).Value = 77
$ws.Cells.Item(475,4).Value = 529
$ws.Cells.Item(475,5).Value = 6
$ws.Cells.Item(475,6).Value = 1
$ws.Cells.Item(476,1).Value = "2024-11-25 23:36:32"
$ws.Cells.Item(476,2).Value = "020917.NC"
$ws.Cells.Item(476,3).Value = 661
$ws.Cells.Item(476,4).Value = 207
$ws.Cells.Item(476,5).Value = 32
$ws.Cells.Item(476,6).Value = 1
$ws.Cells.Item(477,1).Value = "2024-11-26 15:57:14"
$ws.Cells.Item(477,2).Value = "020867.NC"
$ws.Cells.Item(477,3).Value = 16
$ws.Cells.Item(477,4).Value = 133
$ws.Cells.Item(477,5).Value = 6
$ws.Cells.Item(477,6).Value = 1
$ws.Cells.Item(478,1).Value = "2024-11-26 16:20:02"
$ws.Cells.Item(478,2).Value = "021137.NC"
$ws.Cells.Item(478,3).Value = 3
$ws.Cells.Item(478,4).Value = 12
$ws.Cells.Item(478,5).Value = 6
$ws.Cells.Item(478,6).Value = 1
$ws.Cells.Item(479,1).Value = "2024-11-26 16:59:00"
$ws.Cells.Item(479,2).Value = "021059.NC"
$ws.Cells.Item(479,3).Value = 9
$ws.Cells.Item(479,4).Value = 35
$ws.Cells.Item(479,5).Value = 6
$ws.Cells.Item(479,6).Value = 1
$ws.Cells.Item(480,1).Value = "2024-11-26 17:02:41"
$ws.Cells.Item(480,2).Value = "020870.NC"
$ws.Cells.Item(480,3).Value = 30
$ws.Cells.Item(480,4).Value = 122
$ws.Cells.Item(480,5).Value = 10
$ws.Cells.Item(480,6).Value = 1
$ws.Cells.Item(481,1).Value = "2024-11-26 19:49:30"
$ws.Cells.Item(481,2).Value = "020869.NC"
$ws.Cells.Item(481,3).Value = 79
$ws.Cells.Item(481,4).Value = 385
$ws.Cells.Item(481,5).Value = 10
$ws.Cells.Item(481,6).Value = 1
$ws.Cells.Item(482,1).Value = "2024-11-27 02:29:57"
$ws.Cells.Item(482,2).Value = "020909.NC"
$ws.Cells.Item(482,3).Value = 506
$ws.Cells.Item(482,4).Value = 155
$ws.Cells.Item(482,5).Value = 32
$ws.Cells.Item(482,6).Value = 1
$ws.Cells.Item(483,1).Value = "2024-11-28 03:38:50"
$ws.Cells.Item(483,2).Value = "021058.NC"
$ws.Cells.Item(483,3).Value = 24
$ws.Cells.Item(483,4).Value = 122
$ws.Cells.Item(483,5).Value = 6
$ws.Cells.Item(483,6).Value = 1
$ws.Cells.Item(484,1).Value = "2024-11-28 05:27:21"
$ws.Cells.Item(484,2).Value = "021057.NC"
$ws.Cells.Item(484,3).Value = 24
$ws.Cells.Item(484,4).Value = 122
$ws.Cells.Item(484,5).Value = 6
$ws.Cells.Item(484,6).Value = 1
$ws.Cells.Item(485,1).Value = "2024-11-28 09:25:12"
$ws.Cells.Item(485,2).Value = "021054.NC"
$ws.Cells.Item(485,3).Value = 55
$ws.Cells.Item(485,4).Value = 257
$ws.Cells.Item(485,5).Value = 6
$ws.Cells.Item(485,6).Value = 1
$ws.Cells.Item(486,1).Value = "2024-11-28 12:10:00"
$ws.Cells.Item(486,2).Value = "020030.NC"
$ws.Cells.Item(486,3).Value = 41
$ws.Cells.Item(486,4).Value = 195
$ws.Cells.Item(486,5).Value = 6
$ws.Cells.Item(486,6).Value = 1
$ws.Cells.Item(487,1).Value = "2024-11-28 15:43:19"
$ws.Cells.Item(487,2).Value = "021424.NC"
$ws.Cells.Item(487,3).Value = 99
$ws.Cells.Item(487,4).Value = 466
$ws.Cells.Item(487,5).Value = 6
$ws.Cells.Item(487,6).Value = 1
$ws.Cells.Item(488,1).Value = "2024-11-28 19:26:51"
$ws.Cells.Item(488,2).Value = "021425.NC"
$ws.Cells.Item(488,3).Value = 49
$ws.Cells.Item(488,4).Value = 242
$ws.Cells.Item(488,5).Value = 6
$ws.Cells.Item(488,6).Value = 1
$ws.Cells.Item(489,1).Value = "2024-11-28 22:48:20"
$ws.Cells.Item(489,2).Value = "020910.NC"
$ws.Cells.Item(489,3).Value = 673
$ws.Cells.Item(489,4).Value = 177
$ws.Cells.Item(489,5).Value = 32
$ws.Cells.Item(489,6).Value = 1
